$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 31 - this shifts current rows 31..48 down to 32..49,
# carrying their formatting (incl. the "TOTAL" row which becomes row 49) and
# growing the used range to A1:O49.
$ws.Rows("31:31").Insert()

# Populate the freshly inserted row 31 with the new "Linda Boudjemai" booking.
$ws.Range("A31").Value = "Linda Boudjemai"
$ws.Range("B31").Value = "Booking"

# telephone is stored as literal text (it keeps its leading "+"), so force
# text via the classic leading-apostrophe trick, then drop the formatting
# that the apostrophe implies so the cell ends up with no special style.
$ws.Range("C31").Formula = "'+33685678541"
$ws.Range("C31").ClearFormats()

$ws.Range("D31").Value = 45870
$ws.Range("E31").Value = 45880
$ws.Range("F31").Value = 10
$ws.Range("G31").Value = 171.56
$ws.Range("H31").Value = 140.77
$ws.Range("I31").Value = 30.79
$ws.Range("J31").Value = 17.95
$ws.Range("K31").Value = 2025
$ws.Range("L31").Value = 8

# The row this booking was inserted above (row 31) has no confirmed ical
# match, so N31/O31 stay blank - clear the date format/style that Insert
# copied down from the row below.
$ws.Range("N31:O31").ClearFormats()

# The row that used to be "Gregory Blanvillain" (old row 32) had its phone
# number stored as literal text "33687762155.0"; it is now row 33 and the
# number becomes a plain numeric value (no more trailing ".0", no text type).
$ws.Range("C33").Value = 33687762155
